$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 327
$ws1.Range("G3").Value = 70
$ws1.Range("F4").Value = 8434
$ws1.Range("G4").Value = 80
$ws1.Range("F5").Value = 6142
$ws1.Range("G5").Value = 78
$ws1.Range("F10").Value = 320
$ws1.Range("F11").Value = 1113

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G4").Value = 123

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 327
$ws4.Range("G3").Value = 70
$ws4.Range("F4").Value = 8434
$ws4.Range("G4").Value = 80
$ws4.Range("F5").Value = 6142
$ws4.Range("G5").Value = 78
$ws4.Range("F10").Value = 320
$ws4.Range("G13").Value = 123
$ws4.Range("F15").Value = 1113
